$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.190.73"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  +0.06%  "
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.461.39"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.03%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.70"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.45%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.44"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +1.45%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.456.35"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +1.49%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +3.92%  "
$ws.Range("E11").Value = "  +2.62%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.19"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E13").Value = "  -0.24%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.11"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("E15").Value = "  +0.83%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.902.48"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +1.07%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.080.66"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +0.34%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.461.32"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +1.42%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.69"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -1.93%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +4.16%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.07"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -0.96%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.09"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -0.13%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.12"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +2.56%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("E25").Value = "  -0.13%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.19"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.33%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +0.79%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "586.16"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -6.97%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.582.66"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0943"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +0.26%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -0.06%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.38"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -2.40%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("E36").Value = "  -0.13%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -3.03%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.41"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -1.40%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.373"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +0.19%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.81"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +1.81%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.32"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +0.36%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.19"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -0.15%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.71"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.77"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -1.54%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -1.70%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0287"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +21.92%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.65"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("E49").Value = "  -0.87%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.603"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +1.72%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.97"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +2.47%  "
